# ACL-questionnaire-alumni.docx edit:
#   - Merge the "GENDER: (  ) Male (  ) Female" line into one clean run
#     (removing the mid-sentence grammar-check split / proofErr markers).
#   - Merge the "Year of Graduation:____________" line into one clean run
#     the same way.
# The "AGE: ... Program/Degree Completed: ..." paragraph and the
# "Current Employment Status" paragraph are left untouched.

$d = $word.ActiveDocument

$range1 = $d.Content
$range1.Find.Execute("GENDER: (  ) Male (  ) Female", $false, $false, $false, $false, $false, $true, 1, $false, "GENDER: (  ) Male (  ) Female", 2) | Out-Null

$range2 = $d.Content
$range2.Find.Execute("Year of Graduation:____________", $false, $false, $false, $false, $false, $true, 1, $false, "Year of Graduation:____________", 2) | Out-Null
